# cierre de 4 SEPT 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21: payment received - fill in payment date and paid amount
$ws.Range("F21").Value = 44436
$ws.Range("G21").Value = 4859

# Row 22: new credit entry - sale date, client (OBRADOR), amount, payment date, amount paid
$ws.Range("A22").Value = 44438
$ws.Range("D22").Value = "OBRADOR"
$ws.Range("E22").Value = 1360
$ws.Range("F22").Value = 44439
$ws.Range("G22").Value = 1360

# Update the selection to match the author's final cursor position
$null = $ws.Range("G23").Select()
